# Apply updated cryptocurrency price/volume data to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) values are stored as plain text in this sheet
# (t="inlineStr"), never as real numbers. Assigning a numeric-looking
# string straight to .Value lets Excel auto-convert it to a Number, so
# force the cell to Text format first, then restore the default
# "Normal" cell style so no stray number-format style is left on the
# cell (matches the original look: no explicit style index).
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 5).Value = "  +0.94%  "
$dVal = "34.481.85"
Set-TextValue 2 4 $dVal

# Row 3
$ws.Cells.Item(3, 5).Value = "  +0.20%  "
$dVal = "1.794.69"
Set-TextValue 3 4 $dVal

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$ws.Cells.Item(5, 5).Value = "  +0.06%  "
$dVal = "226.83"
Set-TextValue 5 4 $dVal

# Row 6
$ws.Cells.Item(6, 5).Value = "  +1.72%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +1.96%  "
$dVal = "32.54"
Set-TextValue 8 4 $dVal

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.42%  "
$dVal = "0.297"
Set-TextValue 9 4 $dVal

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.48%  "
$dVal = "0.0694"
Set-TextValue 10 4 $dVal

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.28%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +0.31%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "Chainlink"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(13, 5).Value = "  -0.86%  "
$dVal = "11.06"
Set-TextValue 13 4 $dVal

# Row 14
$ws.Cells.Item(14, 2).Value = "WrappedEther"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14, 5).Value = "  -0.40%  "
$dVal = "1.782.81"
Set-TextValue 14 4 $dVal

# Row 15
$ws.Cells.Item(15, 5).Value = "  +2.56%  "
$dVal = "0.638"
Set-TextValue 15 4 $dVal

# Row 16
$ws.Cells.Item(16, 5).Value = "  +1.04%  "
$dVal = "34.466.83"
Set-TextValue 16 4 $dVal

# Row 17
$ws.Cells.Item(17, 5).Value = "  +2.25%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +0.95%  "
$dVal = "68.84"
Set-TextValue 18 4 $dVal

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.64%  "
$dVal = "247.16"
Set-TextValue 19 4 $dVal

# Row 20
$ws.Cells.Item(20, 5).Value = "  +2.50%  "
$dVal = "0.0{0}0798" -f [char]0x2083
Set-TextValue 20 4 $dVal

# Row 21
$ws.Cells.Item(21, 5).Value = "  +3.63%  "
$dVal = "11.23"
Set-TextValue 21 4 $dVal

# Row 22
$ws.Cells.Item(22, 5).Value = "  -0.05%  "
$dVal = "1.00"
Set-TextValue 22 4 $dVal

# Row 23
$ws.Cells.Item(23, 5).Value = "  +1.48%  "
$dVal = "4.17"
Set-TextValue 23 4 $dVal

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.27%  "
$dVal = "2.06"
Set-TextValue 24 4 $dVal

# Row 25
$ws.Cells.Item(25, 5).Value = "  +1.70%  "
$dVal = "163.88"
Set-TextValue 25 4 $dVal

# Row 26
$ws.Cells.Item(26, 5).Value = "  +1.34%  "
$dVal = "7.27"
Set-TextValue 26 4 $dVal

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.04%  "
$dVal = "16.52"
Set-TextValue 27 4 $dVal

# Row 28
$ws.Cells.Item(28, 5).Value = "  +2.34%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.04%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 5).Value = "  +0.13%  "
$dVal = "1.23"
Set-TextValue 30 4 $dVal

# Row 31
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 5).Value = "  +3.62%  "
$dVal = "3.80"
Set-TextValue 31 4 $dVal

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.33%  "
$dVal = "0.0521"
Set-TextValue 32 4 $dVal

# Row 33
$ws.Cells.Item(33, 5).Value = "  +6.80%  "
$dVal = "3.87"
Set-TextValue 33 4 $dVal

# Row 34
$ws.Cells.Item(34, 5).Value = "  +1.25%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -1.42%  "
$dVal = "1.439.93"
Set-TextValue 35 4 $dVal

# Row 36
$ws.Cells.Item(36, 5).Value = "  +7.68%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +2.88%  "
$dVal = "0.668"
Set-TextValue 37 4 $dVal

# Row 38
$ws.Cells.Item(38, 5).Value = "  +2.70%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +0.06%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +5.25%  "
$dVal = "84.62"
Set-TextValue 40 4 $dVal

# Row 41
$ws.Cells.Item(41, 5).Value = "  +1.49%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +1.54%  "
$dVal = "0.936"
Set-TextValue 42 4 $dVal

# Row 43
$ws.Cells.Item(43, 5).Value = "  +2.26%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.36%  "
$dVal = "13.46"
Set-TextValue 44 4 $dVal

# Row 45
$ws.Cells.Item(45, 5).Value = "  +3.28%  "
$dVal = "0.0528"
Set-TextValue 45 4 $dVal

# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.81%  "
$dVal = "6.11"
Set-TextValue 46 4 $dVal

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.33%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.09%  "
$dVal = "1.952.13"
Set-TextValue 48 4 $dVal

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.21%  "
$dVal = "105.87"
Set-TextValue 49 4 $dVal

# Row 50
$ws.Cells.Item(50, 2).Value = "PaxDollar"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(50, 5).Value = "  -0.04%  "
$dVal = "1.00"
Set-TextValue 50 4 $dVal

# Row 51
$ws.Cells.Item(51, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(51, 5).Value = "  -5.05%  "
$dVal = "0.0{0}0129" -f [char]0x2086
Set-TextValue 51 4 $dVal
